$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for refreshed crypto data.
# D-column values are forced as text (leading apostrophe) to preserve exact
# formatting (e.g. trailing zeros like "1.00", "0.0530") instead of Excel
# auto-converting numeric-looking strings to numbers.

$ws.Range("D2").Value = "'57.945.28"
$ws.Range("E2").Value = "  -4.35%  "

$ws.Range("D3").Value = "'2.606.02"
$ws.Range("E3").Value = "  -3.43%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'515.68"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("D6").Value = "'141.91"
$ws.Range("E6").Value = "  -2.31%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  -1.77%  "

$ws.Range("D9").Value = "'6.71"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("E10").Value = "  -2.98%  "

$ws.Range("D11").Value = "'0.336"
$ws.Range("E11").Value = "  -0.88%  "

$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("D13").Value = "'3.064.21"
$ws.Range("E13").Value = "  -3.51%  "

$ws.Range("D14").Value = "'57.948.25"
$ws.Range("E14").Value = "  -4.32%  "

$ws.Range("D15").Value = "'20.63"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").Value = "'2.615.39"
$ws.Range("E17").Value = "  -3.71%  "

$ws.Range("D18").Value = "'4.39"
$ws.Range("E18").Value = "  -2.64%  "

$ws.Range("D19").Value = "'333.61"
$ws.Range("E19").Value = "  -3.39%  "

$ws.Range("D20").Value = "'10.32"
$ws.Range("E20").Value = "  -2.58%  "

$ws.Range("D21").Value = "'6.24"
$ws.Range("E21").Value = "  -3.15%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").Value = "'63.87"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("D25").Value = "'0.166"
$ws.Range("E25").Value = "  -2.42%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("D27").Value = "'7.06"
$ws.Range("E27").Value = "  -3.06%  "

$ws.Range("D28").Value = "'0.0₃0783"
$ws.Range("E28").Value = "  -4.31%  "

$ws.Range("D29").Value = "'6.58"
$ws.Range("E29").Value = "  -3.72%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  -1.88%  "

$ws.Range("D32").Value = "'150.78"
$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").Value = "'18.64"
$ws.Range("E33").Value = "  -2.00%  "

$ws.Range("D34").Value = "'4.08"
$ws.Range("E34").Value = "  -3.96%  "

$ws.Range("D35").Value = "'1.17"
$ws.Range("E35").Value = "  -5.79%  "

$ws.Range("D36").Value = "'0.891"
$ws.Range("E36").Value = "  -5.21%  "

$ws.Range("D37").Value = "'36.45"
$ws.Range("E37").Value = "  -1.80%  "

$ws.Range("D38").Value = "'0.838"
$ws.Range("E38").Value = "  -3.98%  "

$ws.Range("D39").Value = "'1.43"
$ws.Range("E39").Value = "  -6.05%  "

$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D44").Value = "'267.59"
$ws.Range("E44").Value = "  -5.26%  "

$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("D46").Value = "'19.07"
$ws.Range("E46").Value = "  -4.81%  "

$ws.Range("D47").Value = "'0.0530"
$ws.Range("E47").Value = "  -1.43%  "

$ws.Range("D48").Value = "'2.027.26"
$ws.Range("E48").Value = "  -5.56%  "

$ws.Range("E49").Value = "  -2.22%  "

$ws.Range("D50").Value = "'4.60"
$ws.Range("E50").Value = "  -4.42%  "

$ws.Range("D51").Value = "'18.18"
$ws.Range("E51").Value = "  -4.75%  "

# Rows 42/43: Stellar and Mantle swap list positions (re-ranked), each also
# receiving updated Price/Volume values.
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.596"
$ws.Range("E42").Value = "  -2.13%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.0963"
$ws.Range("E43").Value = "  -2.47%  "
